$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.404.58'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.44%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.346.00'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.28%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.08'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '106.08'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.35%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.632'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.72%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.614'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -6.84%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.02'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.47%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0920'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.74%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.39'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.76%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.990'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.13%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.88'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -5.75%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.705.57'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.312.32'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.73%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.355.35'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.52%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -4.33%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.99%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '75.91'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.58'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +6.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '255.22'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -7.64%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.29'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -4.44%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.33'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.74%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.36'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.78%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.75'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.23'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.76%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '175.78'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.23%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '36.20'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.91%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0884'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.50%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.85'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -9.32%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.124'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +16.71%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.90%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.58'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -5.35%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0361'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.88%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.83'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -8.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.66'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -5.70%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.236'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.64%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '70.66'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.87%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -8.24%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.06%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '11.90'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -3.82%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '112.08'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -9.08%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.45'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.88%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.10'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.84%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '84.84'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -10.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.27'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.00%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.04'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.86%  '
